# Apply updated statistics values to the "ოზურგეთი" worksheet and fix the
# active cell selection, matching the author's commit ("files updated and
# bug fixed").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4: "რეგისტრირებული ოჯახი" - updated figures
$ws.Range("E4").Value = 11628
$ws.Range("F4").Value = 10167
$ws.Range("G4").Value = 8066
$ws.Range("H4").Value = 8180
$ws.Range("I4").Value = 7958
$ws.Range("J4").Value = 8629
$ws.Range("K4").Value = 9072

# Row 5: "საარსებო შემწეობის მიმღები ოჯახი" - updated figures
$ws.Range("E5").Value = 1791
$ws.Range("F5").Value = 2133
$ws.Range("G5").Value = 2031
$ws.Range("H5").Value = 2059
$ws.Range("I5").Value = 2281
$ws.Range("J5").Value = 2809
$ws.Range("K5").Value = 3592

# Move/refresh the active selection to I9 (cosmetic state saved with file)
$ws.Range("I9").Select()
